$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the Excel table ("Tabela1") by one row - this grows the table's
# ref/autoFilter range from A1:J70 to A1:J71 and keeps the table metadata
# (filters, columns, style) intact.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Seed formatting for the new row by cloning a representative existing
# data row (one with no border/fill banding) then re-pointing it at the
# "Calibri Light" font used by newer rows in this workbook, so the engine
# derives a fresh set of cell styles (new font + new right-aligned
# number/date formats) instead of reusing the banded row 70 styles.
$ws.Range("A52:J52").Copy()
$ws.Range("A71:J71").PasteSpecial(-4122) | Out-Null
$ws.Range("A71").Locked = $true
$ws.Range("A71:J71").Font.Name = "Calibri Light"

# New day's COVID-19 data (date serial 43971 = 2020-05-20).
$ws.Range("A71").Value = 43971
$ws.Range("B71").Value = 72860
$ws.Range("C71").Value = 909
$ws.Range("D71").Value = 1468
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 21
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 2
$ws.Range("I71").Value = 106
$ws.Range("J71").Value = 1

# Match the author's final selection on the freshly-added row.
$ws.Range("A71:J71").Select() | Out-Null
